# Add a new "Adjective" sheet at the end of the workbook, populate it with
# a short list of descriptive adjectives (used for backstory generation),
# and leave it as the active/selected sheet - matching the target commit
# "Added comment for backstory generation".

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the current last sheet so it lands at the
# end of the tab strip (Worksheets.Add() with no args would insert before
# the active sheet instead).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Adjective"

$adjectives = @(
    "Charming",
    "Brutal",
    "Tough",
    "Sneaky",
    "Small",
    "Muscular",
    "Caring",
    "Apathetic",
    "Mean",
    "Magical",
    "Divine",
    "Merciful"
)

for ($i = 0; $i -lt $adjectives.Count; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $adjectives[$i]
}

# Match the saved selection/active-cell state on the new sheet.
$newSheet.Range("J9").Select()
